$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the empty placeholder textbox (TextBox 12) that duplicated TextBox 14's spot
$s.Shapes.Item("TextBox 12").Delete()

# Reposition shapes: shift the three column header boxes and the divider/answer
# rows upward, and widen + shift the right-hand tag/number boxes to line up with
# the new title textbox position.

$sh = $s.Shapes.Item("TextBox 14")
$sh.Left = 221.06063842773438
$sh.Top = 208.03048706054688

$sh = $s.Shapes.Item("TextBox 17")
$sh.Left = 442.920166015625
$sh.Top = 208.03048706054688

$sh = $s.Shapes.Item("TextBox 23")
$sh.Left = 624.10693359375
$sh.Top = 208.03048706054688

$sh = $s.Shapes.Item("Straight Connector 25")
$sh.Left = 225.58213806152344
$sh.Top = 275.17474365234375

$sh = $s.Shapes.Item("TextBox 26")
$sh.Left = 221.06056213378906
$sh.Top = 284.2145690917969

$sh = $s.Shapes.Item("TextBox 27")
$sh.Left = 440.79742431640625
$sh.Top = 284.2145690917969

$sh = $s.Shapes.Item("TextBox 28")
$sh.Left = 625.306396484375
$sh.Top = 283.6644287109375

$sh = $s.Shapes.Item("TextBox 20")
$sh.Left = 749.088623046875
$sh.Top = 31.890316009521484
$sh.Width = 190.1614227294922
$sh.Height = 21.810945510864258

$sh = $s.Shapes.Item("TextBox 29")
$sh.Left = 749.088623046875
$sh.Top = 68.90126037597656
$sh.Width = 190.1614227294922
$sh.Height = 21.810867309570312

$sh = $s.Shapes.Item("TextBox 30")
$sh.Left = 749.088623046875
$sh.Top = 108.96244812011719
$sh.Width = 190.1614227294922
$sh.Height = 21.810945510864258

$sh = $s.Shapes.Item("TextBox 5")
$sh.Left = 729.4569091796875
$sh.Top = 33.39425277709961

$sh = $s.Shapes.Item("TextBox 8")
$sh.Left = 729.4569091796875
$sh.Top = 70.13835144042969

$sh = $s.Shapes.Item("TextBox 9")
$sh.Left = 729.456787109375
$sh.Top = 111.24252319335938
